$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) retains text formatting so numeric-looking
# strings (e.g. "0.9991", "242.10") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.032.51"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").Value = "1.909.89"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").Value = "0.8464"
$ws.Range("E5").Value = "  +11.10%  "

$ws.Range("D6").Value = "242.10"
$ws.Range("E6").Value = "  +0.81%  "

$ws.Range("D7").Value = "0.9991"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").Value = "0.3252"
$ws.Range("E8").Value = "  +6.42%  "

$ws.Range("D9").Value = "26.81"
$ws.Range("E9").Value = "  +5.49%  "

$ws.Range("D10").Value = "0.07072"
$ws.Range("E10").Value = "  +3.55%  "

$ws.Range("D11").Value = "0.08039"
$ws.Range("E11").Value = "  +0.90%  "

$ws.Range("D12").Value = "0.7544"
$ws.Range("E12").Value = "  +2.02%  "

$ws.Range("D13").Value = "1.906.59"
$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("D14").Value = "5.231"
$ws.Range("E14").Value = "  +1.64%  "

$ws.Range("D15").Value = "93.00"
$ws.Range("E15").Value = "  +2.38%  "

$ws.Range("D16").Value = "14.20"
$ws.Range("E16").Value = "  +2.46%  "

$ws.Range("D17").Value = "30.031.16"
$ws.Range("E17").Value = "  +0.67%  "

$ws.Range("D18").Value = "5.972"
$ws.Range("E18").Value = "  +1.08%  "

$ws.Range("E19").Value = "  +1.07%  "

$ws.Range("D20").Value = "0.000007784"
$ws.Range("E20").Value = "  +1.33%  "

$ws.Range("D21").Value = "2.154.63"
$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("D23").Value = "0.9985"
$ws.Range("E23").Value = "  -0.34%  "

$ws.Range("D24").Value = "7.024"
$ws.Range("E24").Value = "  +1.45%  "

$ws.Range("D25").Value = "0.1627"
$ws.Range("E25").Value = "  +25.63%  "

$ws.Range("D26").Value = "169.57"
$ws.Range("E26").Value = "  +1.90%  "

$ws.Range("D27").Value = "9.283"
$ws.Range("E27").Value = "  +0.88%  "

$ws.Range("E28").Value = "  +1.74%  "

$ws.Range("D29").Value = "2.090"
$ws.Range("E29").Value = "  +3.19%  "

$ws.Range("E30").Value = "  -2.29%  "

$ws.Range("D31").Value = "1.515"
$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("D32").Value = "4.312"
$ws.Range("E32").Value = "  +1.58%  "

$ws.Range("D33").Value = "0.05643"
$ws.Range("E33").Value = "  +7.92%  "

$ws.Range("D34").Value = "4.100"
$ws.Range("E34").Value = "  +0.66%  "

$ws.Range("D35").Value = "1.290"
$ws.Range("E35").Value = "  +3.18%  "

$ws.Range("D36").Value = "0.7376"
$ws.Range("E36").Value = "  +1.84%  "

$ws.Range("D37").Value = "2.720"
$ws.Range("E37").Value = "  +0.16%  "

$ws.Range("D38").Value = "0.01919"
$ws.Range("E38").Value = "  +0.25%  "

$ws.Range("E39").Value = "  +0.62%  "

$ws.Range("D40").Value = "0.4448"
$ws.Range("E40").Value = "  +1.02%  "

$ws.Range("D41").Value = "72.59"
$ws.Range("E41").Value = "  +1.12%  "

$ws.Range("D42").Value = "6.020"
$ws.Range("E42").Value = "  -2.03%  "

$ws.Range("D43").Value = "0.8432"
$ws.Range("E43").Value = "  +1.44%  "

$ws.Range("D44").Value = "1.910"
$ws.Range("E44").Value = "  +1.70%  "

$ws.Range("D45").Value = "0.9992"
$ws.Range("E45").Value = "  -0.20%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "7.629"
$ws.Range("E46").Value = "  +0.37%  "

$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "101.38"
$ws.Range("E47").Value = "  +1.55%  "

$ws.Range("D48").Value = "9.783"
$ws.Range("E48").Value = "  +0.48%  "

$ws.Range("D49").Value = "989.64"
$ws.Range("E49").Value = "  +9.70%  "

$ws.Range("D50").Value = "2.060.94"
$ws.Range("E50").Value = "  +0.20%  "

$ws.Range("D51").Value = "36.38"
$ws.Range("E51").Value = "  +1.10%  "
